$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("G24").Value = 1.75
$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 4.33
$ws.Range("J24").Value = 2.5
$ws.Range("M24").Value = 1.08
$ws.Range("N24").Value = 8
$ws.Range("W24").Value = 2.1
$ws.Range("X24").Value = 1.67
$ws.Range("Z24").Value = 7.5
$ws.Range("AH24").Value = 67
$ws.Range("AK24").Value = 21
$ws.Range("AP24").Value = 3.55

# Row 25
$ws.Range("G25").Value = 2.25
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 3.2
$ws.Range("J25").Value = 3
$ws.Range("M25").Value = 1.06
$ws.Range("N25").Value = 10
$ws.Range("O25").Value = 1.3
$ws.Range("P25").Value = 3.4
$ws.Range("Q25").Value = 2.05
$ws.Range("R25").Value = 1.75
$ws.Range("S25").Value = 3.5
$ws.Range("T25").Value = 1.29
$ws.Range("Z25").Value = 10
$ws.Range("AC25").Value = 19
$ws.Range("AD25").Value = 29
$ws.Range("AN25").Value = 29
$ws.Range("AP25").Value = 3
$ws.Range("AQ25").Value = 1.39
$ws.Range("AR25").Value = 1.58
$ws.Range("AS25").Value = 2.41

# Row 31
$ws.Range("G31").Value = 2.05
$ws.Range("H31").Value = 2.7
$ws.Range("I31").Value = 4.33
$ws.Range("L31").Value = 5.5
$ws.Range("M31").Value = 1.18
$ws.Range("N31").Value = 4.5
$ws.Range("R31").Value = 1.27
$ws.Range("AE31").Value = 4.5
$ws.Range("AG31").Value = 26
$ws.Range("AH31").Value = 126
$ws.Range("AL31").Value = 19
$ws.Range("AN31").Value = 51

# Row 48
$ws.Range("G48").Value = 1.57
$ws.Range("H48").Value = 3.75
$ws.Range("J48").Value = 2.25
$ws.Range("O48").Value = 1.44
$ws.Range("P48").Value = 2.63
$ws.Range("Y48").Value = 5
$ws.Range("AP48").Value = 4
$ws.Range("AQ48").Value = 1.23
$ws.Range("AR48").Value = 1.85
$ws.Range("AS48").Value = 2

# Row 108
$ws.Range("G108").Value = 3.3
$ws.Range("J108").Value = 3.6
$ws.Range("O108").Value = 1.18
$ws.Range("P108").Value = 4.5
$ws.Range("Q108").Value = 1.65
$ws.Range("R108").Value = 2.2
$ws.Range("S108").Value = 2.5
$ws.Range("T108").Value = 1.5
$ws.Range("U108").Value = 1.3
$ws.Range("V108").Value = 3.4
$ws.Range("W108").Value = 1.57
$ws.Range("X108").Value = 2.25
$ws.Range("Y108").Value = 13
$ws.Range("Z108").Value = 19
$ws.Range("AA108").Value = 12
$ws.Range("AC108").Value = 23
$ws.Range("AF108").Value = 7
$ws.Range("AG108").Value = 12
$ws.Range("AH108").Value = 34
$ws.Range("AI108").Value = 126
$ws.Range("AJ108").Value = 10
$ws.Range("AN108").Value = 15
$ws.Range("AO108").Value = 21

# Row 109
$ws.Range("G109").Value = 1.45
$ws.Range("J109").Value = 2
$ws.Range("L109").Value = 7
$ws.Range("M109").Value = 1.05
$ws.Range("N109").Value = 11
$ws.Range("O109").Value = 1.29
$ws.Range("P109").Value = 3.5
$ws.Range("Q109").Value = 1.9
$ws.Range("R109").Value = 1.95
$ws.Range("W109").Value = 2.05
$ws.Range("X109").Value = 1.7
$ws.Range("AB109").Value = 9.5
$ws.Range("AF109").Value = 8.5
$ws.Range("AI109").Value = 451

# Row 112
$ws.Range("G112").Value = 4.3
$ws.Range("H112").Value = 3.75
$ws.Range("I112").Value = 1.7
$ws.Range("J112").Value = 4.5
$ws.Range("K112").Value = 2.22
$ws.Range("L112").Value = 2.22
$ws.Range("O112").Value = 1.23
$ws.Range("P112").Value = 3.35
$ws.Range("Q112").Value = 1.7
$ws.Range("R112").Value = 1.93
$ws.Range("S112").Value = 2.62
$ws.Range("T112").Value = 1.37
$ws.Range("Y112").Value = 13
$ws.Range("Z112").Value = 25
$ws.Range("AA112").Value = 14
$ws.Range("AB112").Value = 70
$ws.Range("AC112").Value = 40
$ws.Range("AE112").Value = 11.75
$ws.Range("AF112").Value = 7.4
$ws.Range("AG112").Value = 15
$ws.Range("AH112").Value = 65
$ws.Range("AJ112").Value = 7.7
$ws.Range("AK112").Value = 8.5
$ws.Range("AM112").Value = 13.5
$ws.Range("AN112").Value = 13

# Row 113
$ws.Range("G113").Value = 1.45
$ws.Range("H113").Value = 4.2
$ws.Range("I113").Value = 6.2
$ws.Range("J113").Value = 1.98
$ws.Range("K113").Value = 2.25
$ws.Range("L113").Value = 5.9
$ws.Range("O113").Value = 1.26
$ws.Range("S113").Value = 2.8
$ws.Range("T113").Value = 1.33
$ws.Range("W113").Value = 1.98
$ws.Range("X113").Value = 1.65
$ws.Range("Y113").Value = 6.4
$ws.Range("Z113").Value = 6.5
$ws.Range("AB113").Value = 9.5
$ws.Range("AC113").Value = 12
$ws.Range("AD113").Value = 30
$ws.Range("AF113").Value = 8.25
$ws.Range("AG113").Value = 20
$ws.Range("AH113").Value = 110
$ws.Range("AJ113").Value = 15.5
$ws.Range("AK113").Value = 37
$ws.Range("AL113").Value = 20
$ws.Range("AM113").Value = 120
$ws.Range("AN113").Value = 70
$ws.Range("AO113").Value = 70

# Row 114
$ws.Range("G114").Value = 1.28
$ws.Range("H114").Value = 5
$ws.Range("I114").Value = 9
$ws.Range("L114").Value = 7.4
$ws.Range("Q114").Value = 1.47
$ws.Range("R114").Value = 2.32
$ws.Range("S114").Value = 2.12
$ws.Range("Y114").Value = 8.25
$ws.Range("Z114").Value = 6.9
$ws.Range("AA114").Value = 8.5
$ws.Range("AB114").Value = 8.25
$ws.Range("AE114").Value = 15.5
$ws.Range("AF114").Value = 10.5
$ws.Range("AI114").Value = 600
$ws.Range("AJ114").Value = 28
$ws.Range("AK114").Value = 70
$ws.Range("AL114").Value = 28
$ws.Range("AM114").Value = 250
$ws.Range("AO114").Value = 80

# Row 121
$ws.Range("G121").Value = 1.38
$ws.Range("H121").Value = 4.5
$ws.Range("I121").Value = 8
$ws.Range("J121").Value = 1.91
$ws.Range("L121").Value = 7.5
$ws.Range("Q121").Value = 1.85
$ws.Range("R121").Value = 2
$ws.Range("S121").Value = 3
$ws.Range("T121").Value = 1.36
$ws.Range("W121").Value = 2.1
$ws.Range("X121").Value = 1.67
$ws.Range("Y121").Value = 6.5
$ws.Range("Z121").Value = 6
$ws.Range("AA121").Value = 9
$ws.Range("AB121").Value = 8.5
$ws.Range("AF121").Value = 9
$ws.Range("AG121").Value = 23
$ws.Range("AH121").Value = 67
$ws.Range("AK121").Value = 41
$ws.Range("AL121").Value = 23
$ws.Range("AM121").Value = 101

# Row 171
$ws.Range("H171").Value = 2.87
$ws.Range("I171").Value = 2.75
$ws.Range("J171").Value = 3.35
$ws.Range("K171").Value = 1.88
$ws.Range("L171").Value = 3.45
$ws.Range("O171").Value = 1.5
$ws.Range("P171").Value = 2.27
$ws.Range("Q171").Value = 2.42
$ws.Range("R171").Value = 1.44
$ws.Range("S171").Value = 4.15
$ws.Range("T171").Value = 1.15
$ws.Range("U171").Value = 1.53
$ws.Range("V171").Value = 2.18
$ws.Range("W171").Value = 2.02
$ws.Range("X171").Value = 1.62
$ws.Range("Y171").Value = 6.5
$ws.Range("Z171").Value = 11.75
$ws.Range("AA171").Value = 10.5
$ws.Range("AB171").Value = 30
$ws.Range("AC171").Value = 28
$ws.Range("AD171").Value = 45
$ws.Range("AE171").Value = 6.3
$ws.Range("AF171").Value = 5.8
$ws.Range("AJ171").Value = 6.5
$ws.Range("AK171").Value = 12
$ws.Range("AM171").Value = 32
